$d = $word.ActiveDocument

# --- Recolor themed heading / title / hyperlink styles to black ---
# (Word keeps the accent1 theme link attributes on the <w:color> element,
#  but the commit forces the cached RGB fallback to pure black.)
$styleNames = @(
    "Title",
    "Heading 1",
    "Heading 2",
    "Heading 3",
    "Heading 4",
    "Heading 5",
    "Heading 6",
    "Heading 7",
    "Heading 8",
    "Heading 9",
    "Hyperlink",
    "TOC Heading"
)

foreach ($name in $styleNames) {
    $style = $d.Styles.Item($name)
    $style.Font.TextColor.RGB = 0
}

# --- Swap the theme's major/minor latin fonts for CMU Serif ---
$fontScheme = $d.DocumentTheme.ThemeFontScheme
$fontScheme.MajorFont.Latin = "CMU Serif"
$fontScheme.MinorFont.Latin = "CMU Serif"
